$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6,8).Value = 44.833332
$ws.Cells.Item(6,9).Value = 44.833332
$ws.Cells.Item(6,11).Value = 134.499996
$ws.Cells.Item(6,13).Value = -22.49999600000001

$ws.Cells.Item(8,8).Value = 586.1429000000001
$ws.Cells.Item(8,9).Value = 184
$ws.Cells.Item(8,10).Value = 2999
$ws.Cells.Item(8,11).Value = 552
$ws.Cells.Item(8,12).Value = 8997
$ws.Cells.Item(8,13).Value = -413
$ws.Cells.Item(8,14).Value = -9275

$ws.Cells.Item(74,8).Value = 4540
$ws.Cells.Item(74,9).Value = 4540
$ws.Cells.Item(74,11).Value = 4540
$ws.Cells.Item(74,13).Value = -3604

$ws.Cells.Item(77,8).Value = 4540
$ws.Cells.Item(77,9).Value = 4540
$ws.Cells.Item(77,11).Value = 22700
$ws.Cells.Item(77,13).Value = -18020

$ws.Cells.Item(129,8).Value = 2586.7144
$ws.Cells.Item(129,9).Value = 2288.8
$ws.Cells.Item(129,10).Value = 2752.2222
$ws.Cells.Item(129,11).Value = 6866.400000000001
$ws.Cells.Item(129,12).Value = 8256.6666
$ws.Cells.Item(129,13).Value = -1866.400000000001
$ws.Cells.Item(129,14).Value = -18256.6666

$ws.Cells.Item(138,8).Value = 3195.5
$ws.Cells.Item(138,10).Value = 3483.389
$ws.Cells.Item(138,12).Value = 10450.167
$ws.Cells.Item(138,14).Value = -20730.167

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(11,8).Value = 30000000
$ws.Cells.Item(11,9).Value = 30000000
$ws.Cells.Item(11,11).Value = 30000000
$ws.Cells.Item(11,13).Value = -29999856

$ws.Cells.Item(13,8).Value = 10000000
$ws.Cells.Item(13,10).Value = 0
$ws.Cells.Item(13,12).Value = 0
$ws.Cells.Item(13,14).ClearContents()

$ws.Cells.Item(74,8).Value = 1721.8235
$ws.Cells.Item(74,9).Value = 1328.7693
$ws.Cells.Item(74,11).Value = 1328.7693
$ws.Cells.Item(74,13).Value = -454.7692999999999

$ws.Cells.Item(77,8).Value = 1721.8235
$ws.Cells.Item(77,9).Value = 1328.7693
$ws.Cells.Item(77,11).Value = 6643.8465
$ws.Cells.Item(77,13).Value = -2275.8465

$ws.Cells.Item(111,8).Value = 24500
$ws.Cells.Item(111,10).Value = 24500
$ws.Cells.Item(111,12).Value = 24500
$ws.Cells.Item(111,14).Value = -32680

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(8,8).Value = 17000
$ws.Cells.Item(8,9).Value = 17000
$ws.Cells.Item(8,11).Value = 17000
$ws.Cells.Item(8,13).Value = -16860

$ws.Cells.Item(12,8).Value = 18000
$ws.Cells.Item(12,9).Value = 0
$ws.Cells.Item(12,11).Value = 0
$ws.Cells.Item(12,13).ClearContents()

$ws.Cells.Item(20,8).Value = 982
$ws.Cells.Item(20,9).Value = 901.5
$ws.Cells.Item(20,10).Value = 1062.5
$ws.Cells.Item(20,11).Value = 901.5
$ws.Cells.Item(20,12).Value = 1062.5
$ws.Cells.Item(20,13).Value = -654.5
$ws.Cells.Item(20,14).Value = -1556.5

$ws.Cells.Item(99,8).Value = 2434.1052
$ws.Cells.Item(99,9).Value = 1633
$ws.Cells.Item(99,11).Value = 1633
$ws.Cells.Item(99,13).Value = -135

$ws.Cells.Item(134,8).Value = 1458.7273
$ws.Cells.Item(134,9).Value = 1458.9524
$ws.Cells.Item(134,11).Value = 4376.857199999999
$ws.Cells.Item(134,13).Value = -1841.857199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(10,8).Value = 11001.5
$ws.Cells.Item(10,9).Value = 9995
$ws.Cells.Item(10,10).Value = 12008
$ws.Cells.Item(10,11).Value = 9995
$ws.Cells.Item(10,12).Value = 12008
$ws.Cells.Item(10,13).Value = -9856
$ws.Cells.Item(10,14).Value = -12286

$ws.Cells.Item(11,8).Value = 9699.666999999999
$ws.Cells.Item(11,9).Value = 10101
$ws.Cells.Item(11,11).Value = 10101
$ws.Cells.Item(11,13).Value = -9961

$ws.Cells.Item(13,8).Value = 9999
$ws.Cells.Item(13,9).Value = 0
$ws.Cells.Item(13,10).Value = 9999
$ws.Cells.Item(13,11).Value = 0
$ws.Cells.Item(13,12).Value = 9999
$ws.Cells.Item(13,13).ClearContents()
$ws.Cells.Item(13,14).Value = -10277

$ws.Cells.Item(31,8).Value = 1448.25
$ws.Cells.Item(31,9).Value = 1431
$ws.Cells.Item(31,10).Value = 1500
$ws.Cells.Item(31,11).Value = 1431
$ws.Cells.Item(31,12).Value = 1500
$ws.Cells.Item(31,13).Value = -1136
$ws.Cells.Item(31,14).Value = -2090

$ws.Cells.Item(34,8).Value = 1448.25
$ws.Cells.Item(34,9).Value = 1431
$ws.Cells.Item(34,10).Value = 1500
$ws.Cells.Item(34,11).Value = 1431
$ws.Cells.Item(34,12).Value = 1500
$ws.Cells.Item(34,13).Value = -1229
$ws.Cells.Item(34,14).Value = -1904

$ws.Cells.Item(58,8).Value = 2951.7144
$ws.Cells.Item(58,9).Value = 2732.4
$ws.Cells.Item(58,11).Value = 2732.4
$ws.Cells.Item(58,13).Value = -2529.4

$ws.Cells.Item(99,8).Value = 1901.1666
$ws.Cells.Item(99,9).Value = 1906.2
$ws.Cells.Item(99,10).Value = 1876
$ws.Cells.Item(99,11).Value = 1906.2
$ws.Cells.Item(99,12).Value = 1876
$ws.Cells.Item(99,13).Value = -408.2
$ws.Cells.Item(99,14).Value = -4872

$ws.Cells.Item(122,8).Value = 4559.125
$ws.Cells.Item(122,9).Value = 5121.75
$ws.Cells.Item(122,11).Value = 15365.25
$ws.Cells.Item(122,13).Value = -12915.25

$ws.Cells.Item(126,8).Value = 1901.1666
$ws.Cells.Item(126,9).Value = 1906.2
$ws.Cells.Item(126,10).Value = 1876
$ws.Cells.Item(126,11).Value = 5718.6
$ws.Cells.Item(126,12).Value = 5628
$ws.Cells.Item(126,13).Value = -3248.6
$ws.Cells.Item(126,14).Value = -10568

$ws.Cells.Item(132,8).Value = 2000
$ws.Cells.Item(132,9).Value = 2000
$ws.Cells.Item(132,11).Value = 6000
$ws.Cells.Item(132,13).Value = -3470

$ws.Cells.Item(136,8).Value = 2951.7144
$ws.Cells.Item(136,9).Value = 2732.4
$ws.Cells.Item(136,11).Value = 8197.200000000001
$ws.Cells.Item(136,13).Value = -5647.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4,8).Value = 1983425.2
$ws.Cells.Item(4,9).Value = 1134438.1
$ws.Cells.Item(4,11).Value = 3403314.3
$ws.Cells.Item(4,13).Value = -3403202.3

$ws.Cells.Item(6,8).Value = 3.6
$ws.Cells.Item(6,9).Value = 3.6
$ws.Cells.Item(6,11).Value = 10.8
$ws.Cells.Item(6,13).Value = 102.2

$ws.Cells.Item(10,8).Value = 468.27274
$ws.Cells.Item(10,9).Value = 127.77778
$ws.Cells.Item(10,10).Value = 2000.5
$ws.Cells.Item(10,11).Value = 383.33334
$ws.Cells.Item(10,12).Value = 6001.5
$ws.Cells.Item(10,13).Value = -244.33334
$ws.Cells.Item(10,14).Value = -6279.5

$ws.Cells.Item(113,8).Value = 946.8
$ws.Cells.Item(113,9).Value = 914.5
$ws.Cells.Item(113,10).Value = 954.875
$ws.Cells.Item(113,11).Value = 2743.5
$ws.Cells.Item(113,12).Value = 2864.625
$ws.Cells.Item(113,13).Value = -573.5
$ws.Cells.Item(113,14).Value = -7204.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2,8).Value = 1195.6
$ws.Cells.Item(2,9).Value = 494.5
$ws.Cells.Item(2,11).Value = 494.5
$ws.Cells.Item(2,13).Value = -381.5

$ws.Cells.Item(3,8).Value = 336334.34
$ws.Cells.Item(3,10).Value = 0
$ws.Cells.Item(3,12).Value = 0
$ws.Cells.Item(3,14).ClearContents()

$ws.Cells.Item(9,8).Value = 1483
$ws.Cells.Item(9,9).Value = 974.6
$ws.Cells.Item(9,11).Value = 974.6
$ws.Cells.Item(9,13).Value = -804.6

$ws.Cells.Item(46,8).Value = 22500
$ws.Cells.Item(46,10).Value = 42000
$ws.Cells.Item(46,12).Value = 42000
$ws.Cells.Item(46,14).Value = -42312

$ws.Cells.Item(122,8).Value = 1798.8572
$ws.Cells.Item(122,9).Value = 1798.8572
$ws.Cells.Item(122,11).Value = 5396.571599999999
$ws.Cells.Item(122,13).Value = -2946.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7,8).Value = 0
$ws.Cells.Item(7,9).Value = 0
$ws.Cells.Item(7,10).Value = 0
$ws.Cells.Item(7,11).Value = 0
$ws.Cells.Item(7,12).Value = 0
$ws.Cells.Item(7,13).ClearContents()
$ws.Cells.Item(7,14).ClearContents()

$ws.Cells.Item(10,8).Value = 6003
$ws.Cells.Item(10,9).Value = 6003
$ws.Cells.Item(10,11).Value = 6003
$ws.Cells.Item(10,13).Value = -5863

$ws.Cells.Item(40,8).Value = 18336
$ws.Cells.Item(40,9).Value = 2504
$ws.Cells.Item(40,11).Value = 2504
$ws.Cells.Item(40,13).Value = -2368

$ws.Cells.Item(68,8).Value = 2160.3333
$ws.Cells.Item(68,9).Value = 2160.3333
$ws.Cells.Item(68,11).Value = 2160.3333
$ws.Cells.Item(68,13).Value = -1411.3333

$ws.Cells.Item(71,8).Value = 2160.3333
$ws.Cells.Item(71,9).Value = 2160.3333
$ws.Cells.Item(71,11).Value = 10801.6665
$ws.Cells.Item(71,13).Value = -7057.666499999999

$ws.Cells.Item(101,8).Value = 22321
$ws.Cells.Item(101,10).Value = 22321
$ws.Cells.Item(101,12).Value = 22321
$ws.Cells.Item(101,14).Value = -28811

$ws.Cells.Item(126,8).Value = 0
$ws.Cells.Item(126,9).Value = 0
$ws.Cells.Item(126,10).Value = 0
$ws.Cells.Item(126,11).Value = 0
$ws.Cells.Item(126,12).Value = 0
$ws.Cells.Item(126,13).ClearContents()
$ws.Cells.Item(126,14).ClearContents()

$ws.Cells.Item(132,8).Value = 2626.96
$ws.Cells.Item(132,9).Value = 2256.1667
$ws.Cells.Item(132,10).Value = 2969.2307
$ws.Cells.Item(132,11).Value = 6768.500100000001
$ws.Cells.Item(132,12).Value = 8907.6921
$ws.Cells.Item(132,13).Value = -4238.500100000001
$ws.Cells.Item(132,14).Value = -13967.6921

$ws.Cells.Item(136,8).Value = 55633944
$ws.Cells.Item(136,9).Value = 67583
$ws.Cells.Item(136,11).Value = 202749
$ws.Cells.Item(136,13).Value = -200199

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(8,8).Value = 0
$ws.Cells.Item(8,9).Value = 0
$ws.Cells.Item(8,11).Value = 0
$ws.Cells.Item(8,13).ClearContents()

$ws.Cells.Item(81,8).Value = 2248.75
$ws.Cells.Item(81,10).Value = 999
$ws.Cells.Item(81,12).Value = 1998
$ws.Cells.Item(81,14).Value = -4120

$ws.Cells.Item(84,8).Value = 2248.75
$ws.Cells.Item(84,10).Value = 999
$ws.Cells.Item(84,12).Value = 9990
$ws.Cells.Item(84,14).Value = -20598

$ws.Cells.Item(115,8).Value = 50000
$ws.Cells.Item(115,10).Value = 50000
$ws.Cells.Item(115,12).Value = 50000
$ws.Cells.Item(115,14).Value = -53134

$ws.Cells.Item(126,8).Value = 1059.3529
$ws.Cells.Item(126,9).Value = 1063.0625
$ws.Cells.Item(126,11).Value = 3189.1875
$ws.Cells.Item(126,13).Value = -719.1875

$ws.Cells.Item(132,8).Value = 4176.467
$ws.Cells.Item(132,9).Value = 3117.2856
$ws.Cells.Item(132,10).Value = 19005
$ws.Cells.Item(132,11).Value = 9351.856800000001
$ws.Cells.Item(132,12).Value = 57015
$ws.Cells.Item(132,13).Value = -6821.856800000001
$ws.Cells.Item(132,14).Value = -62075

$ws.Cells.Item(136,8).Value = 4177.3125
$ws.Cells.Item(136,9).Value = 4411.4614
$ws.Cells.Item(136,11).Value = 13234.3842
$ws.Cells.Item(136,13).Value = -10684.3842
